$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: clear value of I13 (keep its style), was "Account"
$ws.Range("I13").ClearContents()

# Row 14: remove formula from I14 (was =$A14 -> "Acc_Stef"), and instead
# put "Student" shared string value in H14
$ws.Range("I14").Clear()
$ws.Range("H14").Value = "Student"

# Row 15: remove value from G15 (was "Account manager")
$ws.Range("G15").Clear()

# Update the sheet's selected cell to G15
$ws.Range("G15").Select()
